$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 440
$ws.Range("I6").Value = 350
$ws.Range("J6").Value = 800
$ws.Range("K6").Value = 1050
$ws.Range("L6").Value = 2400
$ws.Range("M6").Value = -938
$ws.Range("N6").Value = -2624

$ws.Range("H96").Value = 397.65384
$ws.Range("I96").Value = 320.5625
$ws.Range("J96").Value = 521
$ws.Range("K96").Value = 961.6875
$ws.Range("L96").Value = 1563
$ws.Range("M96").Value = 411.3125
$ws.Range("N96").Value = -4309

$ws.Range("H101").Value = 1068.0435
$ws.Range("I101").Value = 919.41174
$ws.Range("J101").Value = 1489.1666
$ws.Range("K101").Value = 2758.23522
$ws.Range("L101").Value = 4467.4998
$ws.Range("M101").Value = -1136.23522
$ws.Range("N101").Value = -7711.4998

$ws.Range("H123").Value = 45449.355
$ws.Range("J123").Value = 46631
$ws.Range("L123").Value = 46631
$ws.Range("N123").Value = -56431

$ws.Range("H126").Value = 43670
$ws.Range("J126").Value = 43670
$ws.Range("L126").Value = 43670
$ws.Range("N126").Value = -53550

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -884

$ws.Range("H5").Value = 126.666664
$ws.Range("I5").Value = 126.666664
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 126.666664
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -14.666664

$ws.Range("H37").Value = 19343.334
$ws.Range("J37").Value = 20456.111
$ws.Range("L37").Value = 20456.111
$ws.Range("N37").Value = -21002.111

$ws.Range("H74").Value = 1514.8163
$ws.Range("I74").Value = 616.13336
$ws.Range("J74").Value = 11625
$ws.Range("K74").Value = 616.13336
$ws.Range("L74").Value = 11625
$ws.Range("M74").Value = 257.86664
$ws.Range("N74").Value = -13373

$ws.Range("H77").Value = 1514.8163
$ws.Range("I77").Value = 616.13336
$ws.Range("J77").Value = 11625
$ws.Range("K77").Value = 3080.6668
$ws.Range("L77").Value = 58125
$ws.Range("M77").Value = 1287.3332
$ws.Range("N77").Value = -66861

$ws.Range("H110").Value = 1357.3667
$ws.Range("I110").Value = 1462.3462
$ws.Range("J110").Value = 675
$ws.Range("K110").Value = 1462.3462
$ws.Range("L110").Value = 675
$ws.Range("M110").Value = 582.6538
$ws.Range("N110").Value = -4765

$ws.Range("H118").Value = 39498
$ws.Range("J118").Value = 39498
$ws.Range("L118").Value = 39498
$ws.Range("N118").Value = -42812

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 126.666664
$ws.Range("I4").Value = 126.666664
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 126.666664
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -11.666664

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = 0

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("N68").Value = 0

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("N71").Value = 0

$ws.Range("H75").Value = 30000
$ws.Range("I75").Value = 30000
$ws.Range("K75").Value = 30000
$ws.Range("M75").Value = -29064

$ws.Range("H78").Value = 30000
$ws.Range("I78").Value = 30000
$ws.Range("K78").Value = 90000
$ws.Range("M78").Value = -85320

$ws.Range("H82").Value = 22415.875
$ws.Range("I82").Value = 10308.5
$ws.Range("J82").Value = 29680.3
$ws.Range("K82").Value = 10308.5
$ws.Range("L82").Value = 29680.3
$ws.Range("M82").Value = -9925.5
$ws.Range("N82").Value = -30446.3

$ws.Range("H85").Value = 22415.875
$ws.Range("I85").Value = 10308.5
$ws.Range("J85").Value = 29680.3
$ws.Range("K85").Value = 10308.5
$ws.Range("L85").Value = 29680.3
$ws.Range("M85").Value = -8982.5
$ws.Range("N85").Value = -32332.3

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").ClearContents()
$ws.Range("N116").Value = 0

$ws.Range("H126").Value = 32771.11
$ws.Range("J126").Value = 32771.11
$ws.Range("L126").Value = 32771.11
$ws.Range("N126").Value = -42651.11

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 30.5
$ws.Range("I7").Value = 29.6
$ws.Range("J7").Value = 35
$ws.Range("K7").Value = 29.6
$ws.Range("L7").Value = 35
$ws.Range("M7").Value = 83.40000000000001
$ws.Range("N7").Value = -261

$ws.Range("H60").Value = 23218.846
$ws.Range("J60").Value = 25146
$ws.Range("L60").Value = 25146
$ws.Range("N60").Value = -26168

$ws.Range("H109").Value = 12000
$ws.Range("J109").Value = 12000
$ws.Range("L109").Value = 12000
$ws.Range("N109").Value = -14080

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 10660.556
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 10660.556
$ws.Range("K9").Value = 0
$ws.Range("L9").ClearContents()
$ws.Range("M9").Value = 31981.668
$ws.Range("N9").Value = -32429.668

$ws.Range("H131").Value = 917.33685
$ws.Range("I131").Value = 381.42856
$ws.Range("J131").Value = 959.9659
$ws.Range("K131").Value = 1144.28568
$ws.Range("L131").Value = 2879.8977
$ws.Range("M131").Value = 3895.71432
$ws.Range("N131").Value = -12959.8977

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 16146.818
$ws.Range("J57").Value = 18611.666
$ws.Range("L57").Value = 18611.666
$ws.Range("N57").Value = -20251.666

$ws.Range("H97").Value = 2750
$ws.Range("I97").Value = 3125
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 3125
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -2629
$ws.Range("N97").Value = -2992

$ws.Range("H122").Value = 1718.5555
$ws.Range("I122").Value = 1648.04
$ws.Range("J122").Value = 2600
$ws.Range("K122").Value = 4944.12
$ws.Range("L122").Value = 7800
$ws.Range("M122").Value = -2494.12
$ws.Range("N122").Value = -12700

$ws.Range("H123").Value = 38979.5
$ws.Range("J123").Value = 38979.5
$ws.Range("L123").Value = 38979.5
$ws.Range("N123").Value = -43879.5

$ws.Range("H132").Value = 2210.182
$ws.Range("I132").Value = 2053.9443
$ws.Range("J132").Value = 2397.6667
$ws.Range("K132").Value = 6161.8329
$ws.Range("L132").Value = 7193.000100000001
$ws.Range("M132").Value = -3631.8329
$ws.Range("N132").Value = -12253.0001

$ws.Range("H135").Value = 58286.47
$ws.Range("J135").Value = 58286.47
$ws.Range("L135").Value = 58286.47
$ws.Range("N135").Value = -68426.47

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 33976
$ws.Range("J80").Value = 33976
$ws.Range("L80").Value = 33976
$ws.Range("N80").Value = -36222

$ws.Range("H83").Value = 33976
$ws.Range("J83").Value = 33976
$ws.Range("L83").Value = 101928
$ws.Range("N83").Value = -113160

$ws.Range("H132").Value = 4429.5713
$ws.Range("I132").Value = 4140.478
$ws.Range("K132").Value = 12421.434
$ws.Range("M132").Value = -9891.434000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 1000
$ws.Range("I43").Value = 1000
$ws.Range("K43").Value = 1000
$ws.Range("M43").Value = -851

$ws.Range("H109").Value = 29977
$ws.Range("J109").Value = 29977
$ws.Range("L109").Value = 29977
$ws.Range("N109").Value = -32751
